$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Cells.Item(62, 8).Value = 1645.8
$ws.Cells.Item(62, 10).Value = 1835.6
$ws.Cells.Item(62, 12).Value = 1835.6
$ws.Cells.Item(62, 14).Value = -3083.6
# Row 65
$ws.Cells.Item(65, 8).Value = 1645.8
$ws.Cells.Item(65, 10).Value = 1835.6
$ws.Cells.Item(65, 12).Value = 9178
$ws.Cells.Item(65, 14).Value = -15418
# Row 111
$ws.Cells.Item(111, 8).Value = 2035.8462
$ws.Cells.Item(111, 9).Value = 2459.4285
$ws.Cells.Item(111, 10).Value = 1541.6666
$ws.Cells.Item(111, 11).Value = 7378.2855
$ws.Cells.Item(111, 12).Value = 4624.9998
$ws.Cells.Item(111, 13).Value = -4311.2855
$ws.Cells.Item(111, 14).Value = -10758.9998
# Row 129
$ws.Cells.Item(129, 8).Value = 1301.1127
$ws.Cells.Item(129, 10).Value = 1611.9615
$ws.Cells.Item(129, 12).Value = 4835.8845
$ws.Cells.Item(129, 14).Value = -14835.8845
# Row 138
$ws.Cells.Item(138, 8).Value = 1483.41
$ws.Cells.Item(138, 9).Value = 606.23914
$ws.Cells.Item(138, 10).Value = 2230.6296
$ws.Cells.Item(138, 11).Value = 1818.71742
$ws.Cells.Item(138, 12).Value = 6691.888800000001
$ws.Cells.Item(138, 13).Value = 3321.28258
$ws.Cells.Item(138, 14).Value = -16971.8888

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 31
$ws.Cells.Item(31, 8).Value = 12665.777
$ws.Cells.Item(31, 9).Value = 4249
$ws.Cells.Item(31, 10).Value = 80000
$ws.Cells.Item(31, 11).Value = 4249
$ws.Cells.Item(31, 12).Value = 80000
$ws.Cells.Item(31, 13).Value = -3955
$ws.Cells.Item(31, 14).Value = -80588

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 22658
$ws.Cells.Item(102, 9).Value = 10601
$ws.Cells.Item(102, 11).Value = 10601
$ws.Cells.Item(102, 13).Value = -7356

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 47.090908
$ws.Cells.Item(2, 9).Value = 9.083333
$ws.Cells.Item(2, 10).Value = 92.7
$ws.Cells.Item(2, 11).Value = 54.499998
$ws.Cells.Item(2, 12).Value = 556.2
$ws.Cells.Item(2, 13).Value = 58.500002
$ws.Cells.Item(2, 14).Value = -782.2
# Row 3
$ws.Cells.Item(3, 8).Value = 5595.9414
$ws.Cells.Item(3, 9).Value = 4925.5557
$ws.Cells.Item(3, 10).Value = 6350.125
$ws.Cells.Item(3, 11).Value = 14776.6671
$ws.Cells.Item(3, 12).Value = 19050.375
$ws.Cells.Item(3, 13).Value = -14664.6671
$ws.Cells.Item(3, 14).Value = -19274.375
# Row 42
$ws.Cells.Item(42, 8).Value = 3142.8572
$ws.Cells.Item(42, 9).Value = 2000
$ws.Cells.Item(42, 10).Value = 3333.3333
$ws.Cells.Item(42, 11).Value = 6000
$ws.Cells.Item(42, 12).Value = 9999.999899999999
$ws.Cells.Item(42, 13).Value = -5466
$ws.Cells.Item(42, 14).Value = -11067.9999
# Row 109
$ws.Cells.Item(109, 8).Value = 1415.5
$ws.Cells.Item(109, 9).Value = 755.1539
$ws.Cells.Item(109, 10).Value = 10000
$ws.Cells.Item(109, 11).Value = 2265.4617
$ws.Cells.Item(109, 12).Value = 30000
$ws.Cells.Item(109, 13).Value = -1225.4617
$ws.Cells.Item(109, 14).Value = -32080
# Row 110
$ws.Cells.Item(110, 8).Value = 13908.5
$ws.Cells.Item(110, 9).Value = 4000
$ws.Cells.Item(110, 10).Value = 14430
$ws.Cells.Item(110, 11).Value = 12000
$ws.Cells.Item(110, 12).Value = 43290
$ws.Cells.Item(110, 13).Value = -7910
$ws.Cells.Item(110, 14).Value = -51470
# Row 111
$ws.Cells.Item(111, 8).Value = 7740.8335
$ws.Cells.Item(111, 9).Value = 812
$ws.Cells.Item(111, 10).Value = 12690
$ws.Cells.Item(111, 11).Value = 2436
$ws.Cells.Item(111, 12).Value = 38070
$ws.Cells.Item(111, 13).Value = 631
$ws.Cells.Item(111, 14).Value = -44204
# Row 112
$ws.Cells.Item(112, 8).Value = 4525.5454
$ws.Cells.Item(112, 9).Value = 3713.5
$ws.Cells.Item(112, 10).Value = 5500
$ws.Cells.Item(112, 11).Value = 11140.5
$ws.Cells.Item(112, 12).Value = 16500
$ws.Cells.Item(112, 13).Value = -10032.5
$ws.Cells.Item(112, 14).Value = -18716
# Row 115
$ws.Cells.Item(115, 8).Value = 5824.5713
$ws.Cells.Item(115, 9).Value = 4002.5
$ws.Cells.Item(115, 10).Value = 8254
$ws.Cells.Item(115, 11).Value = 12007.5
$ws.Cells.Item(115, 12).Value = 24762
$ws.Cells.Item(115, 13).Value = -10832.5
$ws.Cells.Item(115, 14).Value = -27112
# Row 116
$ws.Cells.Item(116, 8).Value = 4998.6
$ws.Cells.Item(116, 9).Value = 4343
$ws.Cells.Item(116, 10).Value = 5982
$ws.Cells.Item(116, 11).Value = 13029
$ws.Cells.Item(116, 12).Value = 17946
$ws.Cells.Item(116, 13).Value = -9587
$ws.Cells.Item(116, 14).Value = -24830
# Row 118
$ws.Cells.Item(118, 8).Value = 2860.853
$ws.Cells.Item(118, 9).Value = 2228.1667
$ws.Cells.Item(118, 10).Value = 2996.4285
$ws.Cells.Item(118, 11).Value = 6684.500100000001
$ws.Cells.Item(118, 12).Value = 8989.2855
$ws.Cells.Item(118, 13).Value = -5441.500100000001
$ws.Cells.Item(118, 14).Value = -11475.2855
# Row 119
$ws.Cells.Item(119, 8).Value = 4000
$ws.Cells.Item(119, 9).Value = 3000
$ws.Cells.Item(119, 10).Value = 5000
$ws.Cells.Item(119, 11).Value = 9000
$ws.Cells.Item(119, 12).Value = 15000
$ws.Cells.Item(119, 13).Value = -4162
$ws.Cells.Item(119, 14).Value = -24676
# Row 120
$ws.Cells.Item(120, 8).Value = 10000
$ws.Cells.Item(120, 9).Value = 5000
$ws.Cells.Item(120, 10).Value = 12500
$ws.Cells.Item(120, 11).Value = 15000
$ws.Cells.Item(120, 12).Value = 37500
$ws.Cells.Item(120, 13).Value = -10162
$ws.Cells.Item(120, 14).Value = -47176
# Row 133
$ws.Cells.Item(133, 8).Value = 13830.177
$ws.Cells.Item(133, 10).Value = 22501.625
$ws.Cells.Item(133, 12).Value = 67504.875
$ws.Cells.Item(133, 14).Value = -77624.875
# Row 136
$ws.Cells.Item(136, 8).Value = 4105.263
$ws.Cells.Item(136, 9).Value = 975
$ws.Cells.Item(136, 11).Value = 2925
$ws.Cells.Item(136, 13).Value = 2175
# Row 138
$ws.Cells.Item(138, 8).Value = 4629.6875
$ws.Cells.Item(138, 9).Value = 1550.9
$ws.Cells.Item(138, 10).Value = 9761
$ws.Cells.Item(138, 11).Value = 4652.700000000001
$ws.Cells.Item(138, 12).Value = 29283
$ws.Cells.Item(138, 13).Value = 487.2999999999993
$ws.Cells.Item(138, 14).Value = -39563

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 1810.0303
$ws.Cells.Item(102, 9).Value = 1689.6
$ws.Cells.Item(102, 10).Value = 2186.375
$ws.Cells.Item(102, 11).Value = 1689.6
$ws.Cells.Item(102, 12).Value = 2186.375
$ws.Cells.Item(102, 13).Value = -67.59999999999991
$ws.Cells.Item(102, 14).Value = -5430.375
# Row 117
$ws.Cells.Item(117, 8).Value = 55154.75
$ws.Cells.Item(117, 10).Value = 55154.75
$ws.Cells.Item(117, 12).Value = 55154.75
$ws.Cells.Item(117, 14).Value = -62038.75
# Row 134
$ws.Cells.Item(134, 8).Value = 40000
$ws.Cells.Item(134, 10).Value = 40000
$ws.Cells.Item(134, 12).Value = 120000
$ws.Cells.Item(134, 14).Value = -125070

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Cells.Item(122, 8).Value = 3368.3547
$ws.Cells.Item(122, 9).Value = 2138.4614
$ws.Cells.Item(122, 10).Value = 4256.6113
$ws.Cells.Item(122, 11).Value = 6415.3842
$ws.Cells.Item(122, 12).Value = 12769.8339
$ws.Cells.Item(122, 13).Value = -3965.3842
$ws.Cells.Item(122, 14).Value = -17669.8339

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Cells.Item(96, 8).Value = 3578.6667
$ws.Cells.Item(96, 9).Value = 2960
$ws.Cells.Item(96, 10).Value = 4285.7144
$ws.Cells.Item(96, 11).Value = 2960
$ws.Cells.Item(96, 12).Value = 4285.7144
$ws.Cells.Item(96, 13).Value = -1587
$ws.Cells.Item(96, 14).Value = -7031.7144

Write-Output "Applied updates to 165 cells"